# Character Sheet Template Multiclass - crit bonus layout tweak
#
#  - "Class Tree" labels (TextBox 22 / TextBox 25) widened + centred
#    so the longer "special effect" wording fits, instead of a fixed
#    narrow "nowrap" box.
#  - "Notes" caption in the stat box shrunk from 10.5pt to 8pt so it
#    doesn't collide with the widened boxes above it.
#  - the cached "today" field on the master/layouts advanced from
#    23.11.2016 -> 28.11.2016 (re-saved a few days later).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------
# 1) TextBox 22 ("Class Tree") - move/widen + center + wrap text
# ---------------------------------------------------------------
$tb22 = $s.Shapes.Item("TextBox 22")
$tb22.Left   = 26.57503986407874
$tb22.Top    = 359.8207874015748
$tb22.Width  = 173.92496490492127
$tb22.Height = 20.599212598425197
$tb22.TextFrame.WordWrap = -1
$tb22.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# ---------------------------------------------------------------
# 2) TextBox 25 ("Class Tree") - same treatment
# ---------------------------------------------------------------
$tb25 = $s.Shapes.Item("TextBox 25")
$tb25.Left   = 226.39322834645668
$tb25.Top    = 359.8207874015748
$tb25.Width  = 173.92496490492127
$tb25.Height = 20.599212598425197
$tb25.TextFrame.WordWrap = -1
$tb25.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# ---------------------------------------------------------------
# 3) "Notes" caption (3rd paragraph of Rectangle 32) 10.5pt -> 8pt
# ---------------------------------------------------------------
$notesShape = $s.Shapes.Item("Rectangle 32")
$notesRange = $notesShape.TextFrame.TextRange
$notesPara3 = $notesRange.Paragraphs(3, 1)
$notesPara3.Font.Size = 8

# ---------------------------------------------------------------
# 4) Re-cache the "today" date field (datetimeFigureOut) on the
#    slide master and on every layout: 23.11.2016 -> 28.11.2016
# ---------------------------------------------------------------
function Set-DateFieldText {
    param($shapes, [string]$text)

    for ($k = 1; $k -le $shapes.Count; $k++) {
        $shape = $shapes.Item($k)
        $isDatePlaceholder = $false
        try {
            if ($shape.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }

        if ($isDatePlaceholder) {
            $shape.TextFrame.TextRange.Text = $text
        }
    }
}

$master = $p.Designs.Item(1).SlideMaster
Set-DateFieldText $master.Shapes "28.11.2016"

for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    Set-DateFieldText $layout.Shapes "28.11.2016"
}
